$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 109
$ws.Range("I33").Value = 112.166664
$ws.Range("J33").Value = 99.5
$ws.Range("K33").Value = 112.166664
$ws.Range("L33").Value = 99.5
$ws.Range("M33").Value = 116.833336
$ws.Range("N33").Value = -557.5

# Sheet ALC, Row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1547.7778
$ws.Range("J86").Value = 1846.6666
$ws.Range("L86").Value = 1846.6666
$ws.Range("N86").Value = -4092.6666

# Sheet ALC, Row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1547.7778
$ws.Range("J89").Value = 1846.6666
$ws.Range("L89").Value = 9233.333000000001
$ws.Range("N89").Value = -20465.333

# Sheet ALC, Row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2340.4285
$ws.Range("I98").Value = 2569.4443
$ws.Range("J98").Value = 966.3333
$ws.Range("K98").Value = 2569.4443
$ws.Range("L98").Value = 966.3333
$ws.Range("M98").Value = -1071.4443
$ws.Range("N98").Value = -3962.3333

# Sheet ALC, Row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 12718.167
$ws.Range("J116").Value = 5333.2856
$ws.Range("L116").Value = 5333.2856
$ws.Range("N116").Value = -12217.2856

# Sheet ALC, Row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2340.4285
$ws.Range("I122").Value = 2569.4443
$ws.Range("J122").Value = 966.3333
$ws.Range("K122").Value = 7708.3329
$ws.Range("L122").Value = 2898.9999
$ws.Range("M122").Value = -5258.3329
$ws.Range("N122").Value = -7798.9999

# Sheet ALC, Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3299.2
$ws.Range("I138").Value = 2771.72
$ws.Range("J138").Value = 4617.9
$ws.Range("K138").Value = 8315.16
$ws.Range("L138").Value = 13853.7
$ws.Range("M138").Value = -3175.16
$ws.Range("N138").Value = -24133.7

# Sheet ARM, Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 933.1951
$ws.Range("I74").Value = 806.7436
$ws.Range("J74").Value = 3399
$ws.Range("K74").Value = 806.7436
$ws.Range("L74").Value = 3399
$ws.Range("M74").Value = 67.25639999999999
$ws.Range("N74").Value = -5147

# Sheet ARM, Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 933.1951
$ws.Range("I77").Value = 806.7436
$ws.Range("J77").Value = 3399
$ws.Range("K77").Value = 4033.718
$ws.Range("L77").Value = 16995
$ws.Range("M77").Value = 334.2820000000002
$ws.Range("N77").Value = -25731

# Sheet ARM, Row 134
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H134").Value = 39532.668
$ws.Range("J134").Value = 39532.668
$ws.Range("L134").Value = 39532.668
$ws.Range("N134").Value = -49672.668

# Sheet BSM, Row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 332.93103
$ws.Range("I94").Value = 332.93103
$ws.Range("K94").Value = 332.93103
$ws.Range("M94").Value = 118.06897

# Sheet CRP, Row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 735.0909
$ws.Range("I16").Value = 755.4286
$ws.Range("J16").Value = 699.5
$ws.Range("K16").Value = 755.4286
$ws.Range("L16").Value = 699.5
$ws.Range("M16").Value = -468.4286
$ws.Range("N16").Value = -1273.5

# Sheet CRP, Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3950.875
$ws.Range("I31").Value = 3002.25
$ws.Range("J31").Value = 4899.5
$ws.Range("K31").Value = 3002.25
$ws.Range("L31").Value = 4899.5
$ws.Range("M31").Value = -2707.25
$ws.Range("N31").Value = -5489.5

# Sheet CRP, Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3950.875
$ws.Range("I34").Value = 3002.25
$ws.Range("J34").Value = 4899.5
$ws.Range("K34").Value = 3002.25
$ws.Range("L34").Value = 4899.5
$ws.Range("M34").Value = -2800.25
$ws.Range("N34").Value = -5303.5

# Sheet CRP, Row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1124.7778
$ws.Range("I94").Value = 1059.8
$ws.Range("J94").Value = 1206
$ws.Range("K94").Value = 1059.8
$ws.Range("L94").Value = 1206
$ws.Range("M94").Value = -608.8
$ws.Range("N94").Value = -2108

# Sheet CRP, Row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 735.0909
$ws.Range("I113").Value = 755.4286
$ws.Range("J113").Value = 699.5
$ws.Range("K113").Value = 755.4286
$ws.Range("L113").Value = 699.5
$ws.Range("M113").Value = 1414.5714
$ws.Range("N113").Value = -5039.5

# Sheet CUL, Row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 945
$ws.Range("J5").Value = 960
$ws.Range("L5").Value = 2880
$ws.Range("N5").Value = -3104

# Sheet CUL, Row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 912.6667
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 912.6667
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 2738.0001
$ws.Range("N107").Value = -6578.0001
$ws.Range("M107").ClearContents()

# Sheet CUL, Row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 840.9231
$ws.Range("J122").Value = 1043
$ws.Range("L122").Value = 9387
$ws.Range("N122").Value = -14287

# Sheet CUL, Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 753.95
$ws.Range("J131").Value = 780.2308
$ws.Range("L131").Value = 2340.6924
$ws.Range("N131").Value = -12420.6924

# Sheet CUL, Row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# Sheet CUL, Row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 945
$ws.Range("J135").Value = 960
$ws.Range("L135").Value = 8640
$ws.Range("N135").Value = -13710

# Sheet CUL, Row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2566.5557
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Sheet GSM, Row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3009.7273
$ws.Range("I80").Value = 2922.1667
$ws.Range("K80").Value = 2922.1667
$ws.Range("M80").Value = -1924.1667

# Sheet GSM, Row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3009.7273
$ws.Range("I83").Value = 2922.1667
$ws.Range("K83").Value = 14610.8335
$ws.Range("M83").Value = -9618.833500000001

# Sheet GSM, Row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 885.9091
$ws.Range("I97").Value = 787.93335
$ws.Range("J97").Value = 1095.8572
$ws.Range("K97").Value = 787.93335
$ws.Range("L97").Value = 1095.8572
$ws.Range("M97").Value = -291.93335
$ws.Range("N97").Value = -2087.8572

# Sheet GSM, Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5250
$ws.Range("I102").Value = 5999.6665
$ws.Range("K102").Value = 5999.6665
$ws.Range("M102").Value = -4377.6665

# Sheet LTW, Row 25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 57008
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 57008
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 57008
$ws.Range("N25").Value = -57468
$ws.Range("M25").ClearContents()

# Sheet LTW, Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11355
$ws.Range("I122").Value = 9897
$ws.Range("K122").Value = 29691
$ws.Range("M122").Value = -27241

# Sheet WVR, Row 26
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 41942.668
$ws.Range("J26").Value = 41942.668
$ws.Range("L26").Value = 41942.668
$ws.Range("N26").Value = -42528.668

# Sheet WVR, Row 29
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 18995
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

# Sheet WVR, Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3666.3333
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 4249.5
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 12748.5
$ws.Range("M132").Value = -4970
$ws.Range("N132").Value = -17808.5
